$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (data rows 0-23), columns B,C,D,E,G,H,I,K,M,N
# Columns F, J, L, O and column A are unchanged (stay 0 / index values).

# Row 2 (data index 0)
$ws.Cells.Item(2, 2).Value = 10.18193212210974
$ws.Cells.Item(2, 3).Value = 6.146048241211081
$ws.Cells.Item(2, 4).Value = 6.000841204185
$ws.Cells.Item(2, 5).Value = 11.6511982098234
$ws.Cells.Item(2, 7).Value = 39.59665430468949
$ws.Cells.Item(2, 8).Value = 16.76183556689347
$ws.Cells.Item(2, 9).Value = 24.27923766111509
$ws.Cells.Item(2, 11).Value = 10.62873368677392
$ws.Cells.Item(2, 13).Value = 14.5893551430761
$ws.Cells.Item(2, 14).Value = 20.05865695884594

# Row 3 (data index 1)
$ws.Cells.Item(3, 2).Value = 9.921873704013338
$ws.Cells.Item(3, 3).Value = 5.910367986378228
$ws.Cells.Item(3, 4).Value = 5.885185272856104
$ws.Cells.Item(3, 5).Value = 11.43826046631249
$ws.Cells.Item(3, 7).Value = 39.40383858322354
$ws.Cells.Item(3, 8).Value = 16.78164993780516
$ws.Cells.Item(3, 9).Value = 24.30586070439011
$ws.Cells.Item(3, 11).Value = 10.45118847703047
$ws.Cells.Item(3, 13).Value = 14.42721212908878
$ws.Cells.Item(3, 14).Value = 20.11560652831324

# Row 4 (data index 2)
$ws.Cells.Item(4, 2).Value = 9.761508525285656
$ws.Cells.Item(4, 3).Value = 5.76261178955103
$ws.Cells.Item(4, 4).Value = 5.814838911208075
$ws.Cells.Item(4, 5).Value = 11.30940041897275
$ws.Cells.Item(4, 7).Value = 39.29660929103007
$ws.Cells.Item(4, 8).Value = 16.79686780302429
$ws.Cells.Item(4, 9).Value = 24.32704653258796
$ws.Cells.Item(4, 11).Value = 10.3434968062104
$ws.Cells.Item(4, 13).Value = 14.33077356911291
$ws.Cells.Item(4, 14).Value = 20.15244166950302

# Row 5 (data index 3)
$ws.Cells.Item(5, 2).Value = 9.69609444889513
$ws.Cells.Item(5, 3).Value = 5.701739062774982
$ws.Cells.Item(5, 4).Value = 5.786383371747627
$ws.Cells.Item(5, 5).Value = 11.25743868586938
$ws.Cells.Item(5, 7).Value = 39.25574944131914
$ws.Cells.Item(5, 8).Value = 16.80383548319731
$ws.Cells.Item(5, 9).Value = 24.33689433767714
$ws.Cells.Item(5, 11).Value = 10.30000511113286
$ws.Cells.Item(5, 13).Value = 14.29230220528945
$ws.Cells.Item(5, 14).Value = 20.16792264089737

# Row 6 (data index 4)
$ws.Cells.Item(6, 2).Value = 9.685231899764633
$ws.Cells.Item(6, 3).Value = 5.691594561628987
$ws.Cells.Item(6, 4).Value = 5.781672380394311
$ws.Cells.Item(6, 5).Value = 11.24884589243061
$ws.Cells.Item(6, 7).Value = 39.2491368439165
$ws.Cells.Item(6, 8).Value = 16.80503870615923
$ws.Cells.Item(6, 9).Value = 24.33860282696687
$ws.Cells.Item(6, 11).Value = 10.29280889933127
$ws.Cells.Item(6, 13).Value = 14.28596533514219
$ws.Cells.Item(6, 14).Value = 20.17052167218656

# Row 7 (data index 5)
$ws.Cells.Item(7, 2).Value = 9.760626439858083
$ws.Cells.Item(7, 3).Value = 5.761793367594276
$ws.Cells.Item(7, 4).Value = 5.814454237179064
$ws.Cells.Item(7, 5).Value = 11.30869731891894
$ws.Cells.Item(7, 7).Value = 39.29604671716593
$ws.Cells.Item(7, 8).Value = 16.79695867062649
$ws.Cells.Item(7, 9).Value = 24.32717442999714
$ws.Cells.Item(7, 11).Value = 10.34290858797151
$ws.Cells.Item(7, 13).Value = 14.3302513208495
$ws.Cells.Item(7, 14).Value = 20.15264854600935

# Row 8 (data index 6)
$ws.Cells.Item(8, 2).Value = 10.09247089537157
$ws.Cells.Item(8, 3).Value = 6.06547829434265
$ws.Cells.Item(8, 4).Value = 5.960850629123674
$ws.Cells.Item(8, 5).Value = 11.57743175664239
$ws.Cells.Item(8, 7).Value = 39.52787484879511
$ws.Cells.Item(8, 8).Value = 16.76803355105865
$ws.Cells.Item(8, 9).Value = 24.28741155098127
$ws.Cells.Item(8, 11).Value = 10.56727847806058
$ws.Cells.Item(8, 13).Value = 14.53282810507838
$ws.Cells.Item(8, 14).Value = 20.0779058901224

# Row 9 (data index 7)
$ws.Cells.Item(9, 2).Value = 10.7332050695761
$ws.Cells.Item(9, 3).Value = 6.632710181803484
$ws.Cells.Item(9, 4).Value = 6.251337325987666
$ws.Cells.Item(9, 5).Value = 12.11603687502466
$ws.Cells.Item(9, 7).Value = 40.06948568851313
$ws.Cells.Item(9, 8).Value = 16.73556961239398
$ws.Cells.Item(9, 9).Value = 24.24792954154081
$ws.Cells.Item(9, 11).Value = 11.01507322106805
$ws.Cells.Item(9, 13).Value = 14.95280277825623
$ws.Cells.Item(9, 14).Value = 19.94612586143804

# Row 10 (data index 8)
$ws.Cells.Item(10, 2).Value = 11.19214901733115
$ws.Cells.Item(10, 3).Value = 7.027404396632909
$ws.Cells.Item(10, 4).Value = 6.464369043896387
$ws.Cells.Item(10, 5).Value = 12.51442277073976
$ws.Cells.Item(10, 7).Value = 40.51806887013146
$ws.Cells.Item(10, 8).Value = 16.72655763392764
$ws.Cells.Item(10, 9).Value = 24.24250069004583
$ws.Cells.Item(10, 11).Value = 11.34528315896595
$ws.Cells.Item(10, 13).Value = 15.27242502838668
$ws.Cells.Item(10, 14).Value = 19.85827986515451

# Row 11 (data index 9)
$ws.Cells.Item(11, 2).Value = 11.39729015621391
$ws.Cells.Item(11, 3).Value = 7.201346600598148
$ws.Cells.Item(11, 4).Value = 6.560703349032379
$ws.Cells.Item(11, 5).Value = 12.69533871321977
$ws.Cells.Item(11, 7).Value = 40.73255310803404
$ws.Cells.Item(11, 8).Value = 16.72568682732371
$ws.Cells.Item(11, 9).Value = 24.24516672506839
$ws.Cells.Item(11, 11).Value = 11.49504331437694
$ws.Cells.Item(11, 13).Value = 15.41962662231335
$ws.Cells.Item(11, 14).Value = 19.82025580968201

# Row 12 (data index 10)
$ws.Cells.Item(12, 2).Value = 11.47436800794165
$ws.Cells.Item(12, 3).Value = 7.266350447661917
$ws.Cells.Item(12, 4).Value = 6.597061129510092
$ws.Cells.Item(12, 5).Value = 12.76373044602018
$ws.Cells.Item(12, 7).Value = 40.8152163446408
$ws.Cells.Item(12, 8).Value = 16.72582149809927
$ws.Cells.Item(12, 9).Value = 24.24691527936119
$ws.Cells.Item(12, 11).Value = 11.55163049523608
$ws.Cells.Item(12, 13).Value = 15.47557652098896
$ws.Cells.Item(12, 14).Value = 19.80613511306948

# Row 13 (data index 11)
$ws.Cells.Item(13, 2).Value = 11.45779603881412
$ws.Cells.Item(13, 3).Value = 7.25238998989603
$ws.Cells.Item(13, 4).Value = 6.589236857115551
$ws.Cells.Item(13, 5).Value = 12.74900739364532
$ws.Cells.Item(13, 7).Value = 40.79735009532344
$ws.Cells.Item(13, 8).Value = 16.72577183922381
$ws.Cells.Item(13, 9).Value = 24.24650582740032
$ws.Cells.Item(13, 11).Value = 11.53944984164452
$ws.Cells.Item(13, 13).Value = 15.46351828135497
$ws.Cells.Item(13, 14).Value = 19.809163891778

# Row 14 (data index 12)
$ws.Cells.Item(14, 2).Value = 11.40364396103309
$ws.Cells.Item(14, 3).Value = 7.206712133559305
$ws.Cells.Item(14, 4).Value = 6.5636971867679
$ws.Cells.Item(14, 5).Value = 12.70096809439707
$ws.Cells.Item(14, 7).Value = 40.73932527932183
$ws.Cells.Item(14, 8).Value = 16.72568859905033
$ws.Cells.Item(14, 9).Value = 24.2452957681467
$ws.Cells.Item(14, 11).Value = 11.49970158195525
$ws.Cells.Item(14, 13).Value = 15.42422578233125
$ws.Cells.Item(14, 14).Value = 19.81908851853157

# Row 15 (data index 13)
$ws.Cells.Item(15, 2).Value = 11.37039318934077
$ws.Cells.Item(15, 3).Value = 7.178618977802162
$ws.Cells.Item(15, 4).Value = 6.548036392015523
$ws.Cells.Item(15, 5).Value = 12.67152527914046
$ws.Cells.Item(15, 7).Value = 40.70396957576832
$ws.Cells.Item(15, 8).Value = 16.7256980939844
$ws.Cells.Item(15, 9).Value = 24.24465081627574
$ws.Cells.Item(15, 11).Value = 11.47533681685817
$ws.Cells.Item(15, 13).Value = 15.40018353722253
$ws.Cells.Item(15, 14).Value = 19.8252038550187

# Row 16 (data index 14)
$ws.Cells.Item(16, 2).Value = 11.17866225967343
$ws.Cells.Item(16, 3).Value = 7.015918906949039
$ws.Cells.Item(16, 4).Value = 6.458058399829287
$ws.Cells.Item(16, 5).Value = 12.50258695592896
$ws.Cells.Item(16, 7).Value = 40.5042568227723
$ws.Cells.Item(16, 8).Value = 16.72667953107814
$ws.Cells.Item(16, 9).Value = 24.24242985238649
$ws.Cells.Item(16, 11).Value = 11.33548179392534
$ws.Cells.Item(16, 13).Value = 15.26283694237546
$ws.Cells.Item(16, 14).Value = 19.86080377765031

# Row 17 (data index 15)
$ws.Cells.Item(17, 2).Value = 11.06005013080899
$ws.Cells.Item(17, 3).Value = 6.914627950052895
$ws.Cells.Item(17, 4).Value = 6.402684001391115
$ws.Cells.Item(17, 5).Value = 12.39881643093733
$ws.Cells.Item(17, 7).Value = 40.38436996784832
$ws.Cells.Item(17, 8).Value = 16.7281087470066
$ws.Cells.Item(17, 9).Value = 24.24238321209993
$ws.Cells.Item(17, 11).Value = 11.24952665371742
$ws.Cells.Item(17, 13).Value = 15.1790044954087
$ws.Cells.Item(17, 14).Value = 19.8831391139326

# Row 18 (data index 16)
$ws.Cells.Item(18, 2).Value = 10.99148992272749
$ws.Cells.Item(18, 3).Value = 6.855843344808096
$ws.Cells.Item(18, 4).Value = 6.370782226272525
$ws.Cells.Item(18, 5).Value = 12.33910516512561
$ws.Cells.Item(18, 7).Value = 40.31639869866331
$ws.Cells.Item(18, 8).Value = 16.72923470955033
$ws.Cells.Item(18, 9).Value = 24.24283977618996
$ws.Cells.Item(18, 11).Value = 11.20004815398432
$ws.Cells.Item(18, 13).Value = 15.13095914858083
$ws.Cells.Item(18, 14).Value = 19.89616820187055

# Row 19 (data index 17)
$ws.Cells.Item(19, 2).Value = 10.96822123096709
$ws.Cells.Item(19, 3).Value = 6.835851664871563
$ws.Cells.Item(19, 4).Value = 6.359973151221368
$ws.Cells.Item(19, 5).Value = 12.31888588387865
$ws.Cells.Item(19, 7).Value = 40.29355553076705
$ws.Cells.Item(19, 8).Value = 16.72966813116935
$ws.Cells.Item(19, 9).Value = 24.24307736168229
$ws.Cells.Item(19, 11).Value = 11.1832906408747
$ws.Cells.Item(19, 13).Value = 15.11472308397259
$ws.Cells.Item(19, 14).Value = 19.90061096681438

# Row 20 (data index 18)
$ws.Cells.Item(20, 2).Value = 11.07271215495606
$ws.Cells.Item(20, 3).Value = 6.925465298191005
$ws.Cells.Item(20, 4).Value = 6.408584361481353
$ws.Cells.Item(20, 5).Value = 12.40986612334278
$ws.Cells.Item(20, 7).Value = 40.39703063308833
$ws.Cells.Item(20, 8).Value = 16.72792514871307
$ws.Cells.Item(20, 9).Value = 24.24233814296183
$ws.Cells.Item(20, 11).Value = 11.25868123143825
$ws.Cells.Item(20, 13).Value = 15.18791108664118
$ws.Cells.Item(20, 14).Value = 19.88074260536407

# Row 21 (data index 19)
$ws.Cells.Item(21, 2).Value = 11.41956676871298
$ws.Cells.Item(21, 3).Value = 7.220152709012932
$ws.Cells.Item(21, 4).Value = 6.571202418913974
$ws.Cells.Item(21, 5).Value = 12.7150821429369
$ws.Cells.Item(21, 7).Value = 40.75632988031936
$ws.Cells.Item(21, 8).Value = 16.72570044441314
$ws.Cells.Item(21, 9).Value = 24.24563113440351
$ws.Cells.Item(21, 11).Value = 11.51138042146017
$ws.Cells.Item(21, 13).Value = 15.4357617095497
$ws.Cells.Item(21, 14).Value = 19.81616586841373

# Row 22 (data index 20)
$ws.Cells.Item(22, 2).Value = 11.64269868618631
$ws.Cells.Item(22, 3).Value = 7.407688526709124
$ws.Cells.Item(22, 4).Value = 6.67675528188663
$ws.Cells.Item(22, 5).Value = 12.91384404036328
$ws.Cells.Item(22, 7).Value = 40.99953389981314
$ws.Cells.Item(22, 8).Value = 16.72695340028427
$ws.Cells.Item(22, 9).Value = 24.25209057034002
$ws.Cells.Item(22, 11).Value = 11.67578734466885
$ws.Cells.Item(22, 13).Value = 15.59893436797451
$ws.Cells.Item(22, 14).Value = 19.77558254924796

# Row 23 (data index 21)
$ws.Cells.Item(23, 2).Value = 11.52395984935059
$ws.Cells.Item(23, 3).Value = 7.308077323026872
$ws.Cells.Item(23, 4).Value = 6.620498727584832
$ws.Cells.Item(23, 5).Value = 12.80784950906743
$ws.Cells.Item(23, 7).Value = 40.86898370402279
$ws.Cells.Item(23, 8).Value = 16.72603700683268
$ws.Cells.Item(23, 9).Value = 24.24824888775888
$ws.Cells.Item(23, 11).Value = 11.588127265151
$ws.Cells.Item(23, 13).Value = 15.51175426109336
$ws.Cells.Item(23, 14).Value = 19.79709441508837

# Row 24 (data index 22)
$ws.Cells.Item(24, 2).Value = 11.0669887953612
$ws.Cells.Item(24, 3).Value = 6.9205674451273
$ws.Cells.Item(24, 4).Value = 6.405917012201827
$ws.Cells.Item(24, 5).Value = 12.40487071451932
$ws.Cells.Item(24, 7).Value = 40.39130377257924
$ws.Cells.Item(24, 8).Value = 16.72800720560765
$ws.Cells.Item(24, 9).Value = 24.24235701297541
$ws.Cells.Item(24, 11).Value = 11.2545426346701
$ws.Cells.Item(24, 13).Value = 15.18388394232086
$ws.Cells.Item(24, 14).Value = 19.88182548053821

# Row 25 (data index 23)
$ws.Cells.Item(25, 2).Value = 10.56155349093587
$ws.Cells.Item(25, 3).Value = 6.482821293356277
$ws.Cells.Item(25, 4).Value = 6.172647774988279
$ws.Cells.Item(25, 5).Value = 11.96954713694004
$ws.Cells.Item(25, 7).Value = 39.91388082761499
$ws.Cells.Item(25, 8).Value = 16.7417488837174
$ws.Cells.Item(25, 9).Value = 24.25447633635424
$ws.Cells.Item(25, 11).Value = 10.89346754951428
$ws.Cells.Item(25, 13).Value = 14.8370402051784
$ws.Cells.Item(25, 14).Value = 19.98019706543792
